$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text with today's updated rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$currentText = $cellA1.Value2
$updatedText = $currentText.Replace(
    "1000 Bs = 1.9 = 6889.88 pesos",
    "1000 Bs = 1.96 = 7161.1 pesos"
)
$updatedText = $updatedText.Replace(
    "6889.88 pesos = 1.89 = 936.95 Bs",
    "7161.1 pesos = 1.96 = 961.36 Bs"
)
$cellA1.Value2 = $updatedText

# --- Update "tasas" sheet bs/cop <-> usdt rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 509
$wsTasas.Range("O10").Value = 3645
$wsTasas.Range("N12").Value = 3649.98
$wsTasas.Range("O12").Value = 490
